$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1) Date range change: 01/01 - 21/03 Anno 2020  ->  01/01 - 28/03 Anno 2020
# ---------------------------------------------------------------------------
$d.Content.Find.Execute("21/03 Anno 2020", $true, $false, $false, $false, $false, $true, 1, $false, "28/03 Anno 2020", 2)

# ---------------------------------------------------------------------------
# 2) Comuni count: 1.084 -> 1.450
# ---------------------------------------------------------------------------
$d.Content.Find.Execute("1.084", $true, $false, $false, $false, $false, $true, 1, $false, "1.450", 2)

# ---------------------------------------------------------------------------
# 3) Relocate the "_GoBack" bookmark from its old spot (after "I comuni con
#    dati presenti ") to the end of the "20=95-99" list item (after the text
#    changes below have run). Word's own "_GoBack" bookmark simply marks the
#    site of the most recent edit, so we remove it from the old location and
#    recreate it once we reach the new edit location.
# ---------------------------------------------------------------------------
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks("_GoBack").Delete()
}

# ---------------------------------------------------------------------------
# 4) Age-range list: shift every bucket down by one year (5-year buckets
#    became 1 smaller i.e. n=1-5 -> n=1-4, n=6-10 -> n=5-9, etc.)
# ---------------------------------------------------------------------------
$d.Content.Find.Execute("1=1-5", $true, $false, $false, $false, $false, $true, 1, $false, "1=1-4", 2)
$d.Content.Find.Execute("2=6-10", $true, $false, $false, $false, $false, $true, 1, $false, "2=5-9", 2)
$d.Content.Find.Execute("3=11-15", $true, $false, $false, $false, $false, $true, 1, $false, "3=10-14", 2)
$d.Content.Find.Execute("4=16-20", $true, $false, $false, $false, $false, $true, 1, $false, "4=15-19", 2)
$d.Content.Find.Execute("5=21-25", $true, $false, $false, $false, $false, $true, 1, $false, "5=20-24", 2)
$d.Content.Find.Execute("6=26-30", $true, $false, $false, $false, $false, $true, 1, $false, "6=25-29", 2)
$d.Content.Find.Execute("7=31-35", $true, $false, $false, $false, $false, $true, 1, $false, "7=30-34", 2)
$d.Content.Find.Execute("8=36-40", $true, $false, $false, $false, $false, $true, 1, $false, "8=35-39", 2)
$d.Content.Find.Execute("9=41-45", $true, $false, $false, $false, $false, $true, 1, $false, "9=40-44", 2)
$d.Content.Find.Execute("10=46-50", $true, $false, $false, $false, $false, $true, 1, $false, "10=45-49", 2)
$d.Content.Find.Execute("11=51-55", $true, $false, $false, $false, $false, $true, 1, $false, "11=50-54", 2)
$d.Content.Find.Execute("12=56-60", $true, $false, $false, $false, $false, $true, 1, $false, "12=55-59", 2)
$d.Content.Find.Execute("13=61-65", $true, $false, $false, $false, $false, $true, 1, $false, "13=60-64", 2)
$d.Content.Find.Execute("14=66-70", $true, $false, $false, $false, $false, $true, 1, $false, "14=65-69", 2)
$d.Content.Find.Execute("15=71-75", $true, $false, $false, $false, $false, $true, 1, $false, "15=70-74", 2)
$d.Content.Find.Execute("16=76-80", $true, $false, $false, $false, $false, $true, 1, $false, "16=75-79", 2)
$d.Content.Find.Execute("17=81-85", $true, $false, $false, $false, $false, $true, 1, $false, "17=80-84", 2)
$d.Content.Find.Execute("18=86-90", $true, $false, $false, $false, $false, $true, 1, $false, "18=85-89", 2)
$d.Content.Find.Execute("19=91-95", $true, $false, $false, $false, $false, $true, 1, $false, "19=90-94", 2)
$d.Content.Find.Execute("20=96-100", $true, $false, $false, $false, $false, $true, 1, $false, "20=95-99", 2)
$d.Content.Find.Execute("21=101+", $true, $false, $false, $false, $false, $true, 1, $false, "21=100+", 2)

# ---------------------------------------------------------------------------
# 5) Re-insert the "_GoBack" bookmark right after "20=95-99" (mirrors how
#    Word leaves the marker at the spot of the last keystroke). A bookmark
#    collapsed exactly at a paragraph end is placed one character early and
#    then "hopped" over the final character by deleting + retyping that
#    character, which leaves the bookmark collapsed immediately after it.
# ---------------------------------------------------------------------------
$rng = $d.Content
$rng.Find.Execute("20=95-99") | Out-Null
$lastCharStart = $rng.End - 1
$bmTarget = $d.Range($lastCharStart, $lastCharStart)
$d.Bookmarks.Add("_GoBack", $bmTarget)

$lastChar = $d.Range($lastCharStart, $rng.End)
$charText = $lastChar.Text
$lastChar.Delete()
$reinsert = $d.Range($lastCharStart, $lastCharStart)
$reinsert.InsertBefore($charText)

Write-Host "done"
